$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Add a new column G mirroring column F's formatting (header style, data style,
# and totals-row style), then overwrite the header text and values.
$ws.Range("F1:F52").Copy($ws.Range("G1"))

# Header cell
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows (2-51) and totals row (52) are all zero, same as the source column.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Match the new column's width (17 in the saved file) to the rest of the sheet.
$ws.Columns.Item(7).ColumnWidth = 16.1666666666667
